$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text updates (volume/issue number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- Column H width adjustment ---
$ws.Columns.Item(8).ColumnWidth = 7.433768

# --- Row 14 ---
$ws.Range("N14").Value = -50

# --- Row 15 ---
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -21.428571428571

# --- Row 16 ---
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 150
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 27.272727272727
$ws.Range("I16").Value = 108
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = 1.886792452830
$ws.Range("L16").Value = -6.896551724137
$ws.Range("M16").Value = -3.571428571428
$ws.Range("N16").Value = -83.045525902668

# --- Row 17 ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -61.538461538461
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -13.636363636363
$ws.Range("I17").Value = 206
$ws.Range("J17").Value = 199
$ws.Range("K17").Value = 3.517587939698
$ws.Range("L17").Value = 21.176470588235
$ws.Range("M17").Value = 178.378378378378
$ws.Range("N17").Value = -13.080168776371

# --- Row 18 ---
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 113
$ws.Range("J18").Value = 118
$ws.Range("K18").Value = -4.237288135593
$ws.Range("L18").Value = -31.515151515151
$ws.Range("M18").Value = -40.526315789473
$ws.Range("N18").Value = -92.446524064171

# --- Row 19 ---
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 40.476190476190
$ws.Range("I19").Value = 432
$ws.Range("J19").Value = 449
$ws.Range("K19").Value = -3.786191536748
$ws.Range("L19").Value = -10.927835051546
$ws.Range("M19").Value = 57.090909090909
$ws.Range("N19").Value = -59.013282732447

# --- Row 20 ---
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = -32
$ws.Range("I20").Value = 167
$ws.Range("J20").Value = 181
$ws.Range("K20").Value = -7.734806629834
$ws.Range("L20").Value = 9.150326797385
$ws.Range("M20").Value = 19.285714285714
$ws.Range("N20").Value = -93.219650832318

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -21.428571428571
$ws.Range("F21").Value = 120
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = 5.263157894736
$ws.Range("I21").Value = 1039
$ws.Range("J21").Value = 1068
$ws.Range("K21").Value = -2.715355805243
$ws.Range("L21").Value = -5.972850678733
$ws.Range("M21").Value = 30.363864491844
$ws.Range("N21").Value = -82.413676371022

# --- Row 22 ---
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0

# --- Row 23 ---
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -55.555555555555
$ws.Range("I23").Value = 46
$ws.Range("J23").Value = 56
$ws.Range("K23").Value = -17.857142857142
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 142.105263157895

# --- Row 24 ---
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 56.25
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = 9.677419354838
$ws.Range("I24").Value = 812
$ws.Range("J24").Value = 850
$ws.Range("K24").Value = -4.470588235294
$ws.Range("L24").Value = -10.375275938189
$ws.Range("M24").Value = 4.504504504504

# --- Row 25 ---
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = -30.188679245283
$ws.Range("I25").Value = 304
$ws.Range("J25").Value = 363
$ws.Range("K25").Value = -16.253443526170
$ws.Range("L25").Value = -18.059299191374

# --- Row 26 ---
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 44.444444444444
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 41
$ws.Range("I26").Value = 378
$ws.Range("J26").Value = 338
$ws.Range("K26").Value = 11.834319526627
$ws.Range("L26").Value = 32.167832167832
$ws.Range("M26").Value = 3.561643835616

# --- Row 27 ---
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = -36.842105263157

# --- Row 28 ---
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 900
$ws.Range("I28").Value = 42
$ws.Range("K28").Value = 68
$ws.Range("L28").Value = 35.483870967741

# --- Row 29 ---
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -50
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 3
$ws.Range("J29").Value = 6
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = -89.285714285714

# --- Row 30 ---
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 3
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = -25
$ws.Range("L30").Value = -25
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = -85.714285714285

# --- Row 31 ---
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100
$ws.Range("I31").Value = 13
$ws.Range("J31").Value = 23
$ws.Range("K31").Value = -43.478260869565
$ws.Range("L31").Value = 30
